$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 220
$ws.Range("I9").Value = 246.64285
$ws.Range("K9").Value = 246.64285
$ws.Range("M9").Value = -77.64285000000001

$ws.Range("H55").Value = 908.86664
$ws.Range("I55").Value = 774
$ws.Range("J55").Value = 998.7778
$ws.Range("K55").Value = 774
$ws.Range("L55").Value = 998.7778
$ws.Range("M55").Value = -560
$ws.Range("N55").Value = -1426.7778

$ws.Range("H88").Value = 6142.125
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 6142.125
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 6142.125
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -6954.125

$ws.Range("H91").Value = 6142.125
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 6142.125
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 6142.125
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -8950.125

$ws.Range("H129").Value = 1184.7333
$ws.Range("J129").Value = 1250
$ws.Range("L129").Value = 3750
$ws.Range("N129").Value = -13750

$ws.Range("H137").Value = 2972.5
$ws.Range("J137").Value = 5071
$ws.Range("L137").Value = 15213
$ws.Range("N137").Value = -20313

$ws.Range("H138").Value = 2367.5186
$ws.Range("I138").Value = 1866.75
$ws.Range("J138").Value = 3369.0557
$ws.Range("K138").Value = 5600.25
$ws.Range("L138").Value = 10107.1671
$ws.Range("M138").Value = -460.25
$ws.Range("N138").Value = -20387.1671


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 9750.75
$ws.Range("I13").Value = 1998
$ws.Range("J13").Value = 10858.286
$ws.Range("K13").Value = 1998
$ws.Range("L13").Value = 10858.286
$ws.Range("M13").Value = -1854
$ws.Range("N13").Value = -11146.286

$ws.Range("H61").Value = 3998.6667
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H132").Value = 29513.395
$ws.Range("I132").Value = 31557.543
$ws.Range("J132").Value = 5665
$ws.Range("K132").Value = 94672.629
$ws.Range("L132").Value = 16995
$ws.Range("M132").Value = -92142.629
$ws.Range("N132").Value = -22055

$ws.Range("H136").Value = 3998.6667
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3536.6875
$ws.Range("I86").Value = 2494.5
$ws.Range("J86").Value = 4578.875
$ws.Range("K86").Value = 2494.5
$ws.Range("L86").Value = 4578.875
$ws.Range("M86").Value = -1371.5
$ws.Range("N86").Value = -6824.875

$ws.Range("H89").Value = 3536.6875
$ws.Range("I89").Value = 2494.5
$ws.Range("J89").Value = 4578.875
$ws.Range("K89").Value = 12472.5
$ws.Range("L89").Value = 22894.375
$ws.Range("M89").Value = -6856.5
$ws.Range("N89").Value = -34126.375

$ws.Range("H134").Value = 5572.25
$ws.Range("I134").Value = 5368.5713
$ws.Range("J134").Value = 6998
$ws.Range("K134").Value = 16105.7139
$ws.Range("L134").Value = 20994
$ws.Range("M134").Value = -13570.7139
$ws.Range("N134").Value = -26064


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 292.3889
$ws.Range("I7").Value = 193.7
$ws.Range("J7").Value = 415.75
$ws.Range("K7").Value = 193.7
$ws.Range("L7").Value = 415.75
$ws.Range("M7").Value = -80.69999999999999
$ws.Range("N7").Value = -641.75

$ws.Range("H31").Value = 3285.1943
$ws.Range("I31").Value = 2155.0588
$ws.Range("K31").Value = 2155.0588
$ws.Range("M31").Value = -1860.0588

$ws.Range("H34").Value = 3285.1943
$ws.Range("I34").Value = 2155.0588
$ws.Range("K34").Value = 2155.0588
$ws.Range("M34").Value = -1953.0588

$ws.Range("H86").Value = 5122.3335
$ws.Range("I86").Value = 5784
$ws.Range("J86").Value = 3799
$ws.Range("K86").Value = 5784
$ws.Range("L86").Value = 3799
$ws.Range("M86").Value = -4661
$ws.Range("N86").Value = -6045

$ws.Range("H89").Value = 5122.3335
$ws.Range("I89").Value = 5784
$ws.Range("J89").Value = 3799
$ws.Range("K89").Value = 28920
$ws.Range("L89").Value = 18995
$ws.Range("M89").Value = -23304
$ws.Range("N89").Value = -30227

$ws.Range("H107").Value = 2283.8386
$ws.Range("I107").Value = 279.41666
$ws.Range("K107").Value = 279.41666
$ws.Range("M107").Value = 1640.58334


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 9077.923000000001
$ws.Range("J105").Value = 9077.923000000001
$ws.Range("L105").Value = 27233.769
$ws.Range("N105").Value = -32475.769

$ws.Range("H122").Value = 2171.4285
$ws.Range("J122").Value = 2233.5
$ws.Range("L122").Value = 20101.5
$ws.Range("N122").Value = -25001.5

$ws.Range("H127").Value = 8993.666999999999
$ws.Range("J127").Value = 8993.666999999999
$ws.Range("L127").Value = 26981.001
$ws.Range("N127").Value = -36901.001

$ws.Range("H131").Value = 4908.7
$ws.Range("I131").Value = 1196.125
$ws.Range("J131").Value = 7383.75
$ws.Range("K131").Value = 3588.375
$ws.Range("L131").Value = 22151.25
$ws.Range("M131").Value = 1451.625
$ws.Range("N131").Value = -32231.25

$ws.Range("H132").Value = 1436.125
$ws.Range("J132").Value = 1441.2858
$ws.Range("L132").Value = 12971.5722
$ws.Range("N132").Value = -18031.5722

$ws.Range("H134").Value = 861.5833
$ws.Range("J134").Value = 5000
$ws.Range("L134").Value = 15000
$ws.Range("N134").Value = -25140

$ws.Range("H140").Value = 2095.4
$ws.Range("I140").Value = 1799.5652
$ws.Range("J140").Value = 5497.5
$ws.Range("K140").Value = 5398.6956
$ws.Range("L140").Value = 16492.5
$ws.Range("M140").Value = -218.6956
$ws.Range("N140").Value = -26852.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 8125
$ws.Range("I22").Value = 6500
$ws.Range("J22").Value = 9750
$ws.Range("K22").Value = 6500
$ws.Range("L22").Value = 9750
$ws.Range("M22").Value = -5971
$ws.Range("N22").Value = -10808


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 67568.3
$ws.Range("I22").Value = 111706.7
$ws.Range("K22").Value = 111706.7
$ws.Range("M22").Value = -111411.7

$ws.Range("H27").Value = 67568.3
$ws.Range("I27").Value = 111706.7
$ws.Range("K27").Value = 111706.7
$ws.Range("M27").Value = -111599.7


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H7").Value = 999.75
$ws.Range("I7").Value = 999.75
$ws.Range("K7").Value = 999.75
$ws.Range("M7").Value = -886.75

$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -100120

$ws.Range("H136").Value = 5067.75
$ws.Range("I136").Value = 4577.4287
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 13732.2861
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -11182.2861
$ws.Range("N136").Value = -30600

